$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 480, pushing the existing rows
# (480:580) down to (482:582).
$ws.Rows("480:481").Insert()

# Populate the two newly-inserted rows with the new "Femacal de La
# Calera" Coliflor price entries (date 44637), matching the constant
# columns used throughout this sheet.

# Row 480 - Primera
$ws.Cells.Item(480, 1).Value = 3
$ws.Cells.Item(480, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(480, 3).Value = "Coquimbo"
$ws.Cells.Item(480, 4).Value = 44637
$ws.Cells.Item(480, 5).Value = 5
$ws.Cells.Item(480, 6).Value = 100112008
$ws.Cells.Item(480, 7).Value = "Coliflor"
$ws.Cells.Item(480, 8).Value = "Sin especificar"
$ws.Cells.Item(480, 9).Value = "Primera"
$ws.Cells.Item(480, 10).Value = 750
$ws.Cells.Item(480, 11).Value = 1200
$ws.Cells.Item(480, 12).Value = 1200
$ws.Cells.Item(480, 13).Value = 1200
$ws.Cells.Item(480, 14).Value = "$/unidad"
$ws.Cells.Item(480, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(480, 16).Value = 1200
$ws.Cells.Item(480, 17).Value = 1
$ws.Cells.Item(480, 18).Value = "Hortaliza"

# Row 481 - Segunda
$ws.Cells.Item(481, 1).Value = 3
$ws.Cells.Item(481, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(481, 3).Value = "Coquimbo"
$ws.Cells.Item(481, 4).Value = 44637
$ws.Cells.Item(481, 5).Value = 5
$ws.Cells.Item(481, 6).Value = 100112008
$ws.Cells.Item(481, 7).Value = "Coliflor"
$ws.Cells.Item(481, 8).Value = "Sin especificar"
$ws.Cells.Item(481, 9).Value = "Segunda"
$ws.Cells.Item(481, 10).Value = 650
$ws.Cells.Item(481, 11).Value = 900
$ws.Cells.Item(481, 12).Value = 900
$ws.Cells.Item(481, 13).Value = 900
$ws.Cells.Item(481, 14).Value = "$/unidad"
$ws.Cells.Item(481, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(481, 16).Value = 900
$ws.Cells.Item(481, 17).Value = 1
$ws.Cells.Item(481, 18).Value = "Hortaliza"
